# Refresh the cryptocurrency "Price" (D) and "Volume(1h)" (E) columns with
# the latest scraped snapshot.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '96.411.60'
$ws.Range("E2").Value = '  +0.02%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.683.86'
$ws.Range("E3").Value = '  -0.50%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '2.43'
$ws.Range("E4").Value = '  +29.50%  '
$ws.Range("E5").Value = '  -0.05%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '228.43'
$ws.Range("E6").Value = '  -3.17%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '650.49'
$ws.Range("E7").Value = '  +0.02%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.437'
$ws.Range("E8").Value = '  +2.59%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '1.15'
$ws.Range("E9").Value = '  +8.77%  '
$ws.Range("E10").Value = '  -0.06%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '3.680.74'
$ws.Range("E11").Value = '  -0.50%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '47.97'
$ws.Range("E12").Value = '  +8.65%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.210'
$ws.Range("E13").Value = '  +2.43%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.0000301'
$ws.Range("E14").Value = '  -1.98%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '6.66'
$ws.Range("E15").Value = '  -0.84%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '4.387.48'
$ws.Range("E16").Value = '  -0.11%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '96.001.05'
$ws.Range("E17").Value = '  -0.21%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '8.91'
$ws.Range("E18").Value = '  +1.23%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '3.679.84'
$ws.Range("E19").Value = '  -0.50%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '19.62'
$ws.Range("E20").Value = '  +5.30%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '12.99'
$ws.Range("E21").Value = '  -0.05%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.541'
$ws.Range("E22").Value = '  +7.62%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '530.38'
$ws.Range("E23").Value = '  +1.97%  '
$ws.Range("E24").Value = '  -2.19%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.251'
$ws.Range("E25").Value = '  +45.71%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '122.34'
$ws.Range("E26").Value = '  +20.72%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '0.0000209'
$ws.Range("E27").Value = '  +0.12%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '6.86'
$ws.Range("E28").Value = '  -0.78%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '3.880.83'
$ws.Range("E29").Value = '  -0.56%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '13.06'
$ws.Range("E30").Value = '  -0.86%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '13.35'
$ws.Range("E31").Value = '  +10.24%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '3.01'
$ws.Range("E32").Value = '  +0.29%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.999'
$ws.Range("E33").Value = '  -0.16%  '
$ws.Range("E34").Value = '  +0.02%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.84'
$ws.Range("E35").Value = '  -0.47%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '32.97'
$ws.Range("E36").Value = '  +2.60%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '1.00'
$ws.Range("E37").Value = '  -0.03%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.615'
$ws.Range("E38").Value = '  +5.10%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '611.74'
$ws.Range("E39").Value = '  -5.60%  '
$ws.Range("E40").Value = '  +0.00%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '8.52'
$ws.Range("E41").Value = '  -3.36%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '7.15'
$ws.Range("E42").Value = '  +4.42%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.508'
$ws.Range("E43").Value = '  +18.59%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.164'
$ws.Range("E44").Value = '  +2.77%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.0503'
$ws.Range("E45").Value = '  +12.03%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '40.41'
$ws.Range("E46").Value = '  +0.08%  '
$ws.Range("E47").Value = '  -2.02%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.964'
$ws.Range("E48").Value = '  +0.66%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '9.05'
$ws.Range("E49").Value = '  +6.92%  '
$ws.Range("E50").Value = '  +1.11%  '
$ws.Range("E51").Value = '  -0.28%  '
